# Insert a new price record at row 156 (a new weekly observation), shifting
# every subsequent row down by one. The newly inserted row keeps the same
# data as the row immediately below it (which is the former row 156),
# except for a new Fecha (date) value.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a blank row before row 156; everything from 156 down shifts to 157+.
$ws.Range("A156").EntireRow.Insert()

# Populate the newly-inserted row 156 with the same values as the row right
# below it (row 157, which now holds what used to be row 156's data).
$ws.Range("A157:R157").Copy($ws.Range("A156:R156"))

# Set the new row's date (Fecha) to the new observation's value.
$ws.Range("D156").Value = 44741
